$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-07-03 Thursday" "2025-07-04 Friday"

Replace-Text "154×2=308" "250×6=1500"
Replace-Text "746×5=3730" "772×6=4632"
Replace-Text "782×8=6256" "178×4=712"
Replace-Text "576×5=2880" "662×3=1986"
Replace-Text "912×8=7296" "452×5=2260"

Replace-Text "607×4=2428" "588×9=5292"
Replace-Text "448×5=2240" "607×4=2428"
Replace-Text "640×7=4480" "678×8=5424"
Replace-Text "922×5=4610" "443×7=3101"
Replace-Text "967×8=7736" "821×8=6568"

Replace-Text "961×5=4805" "735×5=3675"
Replace-Text "843×3=2529" "350×9=3150"
Replace-Text "278×3=834" "477×9=4293"
Replace-Text "351×4=1404" "312×7=2184"
Replace-Text "627×3=1881" "320×5=1600"

Replace-Text "186×6=1116" "937×7=6559"
Replace-Text "748×9=6732" "140×4=560"
Replace-Text "581×9=5229" "952×4=3808"
Replace-Text "401×2=802" "792×4=3168"
Replace-Text "746×2=1492" "772×2=1544"

Replace-Text "178×6=1068" "417×2=834"
Replace-Text "243×4=972" "212×9=1908"
Replace-Text "597×7=4179" "494×8=3952"
Replace-Text "107×5=535" "244×6=1464"
Replace-Text "263×7=1841" "629×5=3145"
